# update all documents 2024-07-11
#
# Normalizes the product catalogue:
#  - collapses the per-product "subcategoria" labels "Cebolla"/"Papa" into
#    the broader category "Verdura", and "Platano" into "Fruta"
#  - strips the stray leading "/src/" from every path_image value so the
#    paths are relative to the repo root (database_RegistraBOT/...)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row -> old "subcategoria_producto" (column C) replacement -------------
$categoryFixes = @{
    2 = "Verdura"  # GR-VR-BS-GRA-000-001 - Cebolla Morada (was "Cebolla")
    3 = "Verdura"  # GR-VR-BS-GRA-000-002 - Cebolla Blanca (was "Cebolla")
    4 = "Fruta"    # GR-FR-BS-GRA-000-001 - Plátano Seda   (was "Platano")
    5 = "Verdura"  # GR-VR-BS-GRA-000-003 - Papa Amarilla  (was "Papa")
    6 = "Verdura"  # GR-VR-BS-GRA-000-004 - Papa Blanca    (was "Papa")
    7 = "Verdura"  # GR-VR-BS-GRA-000-005 - Papa Huayro    (was "Papa")
    8 = "Verdura"  # GR-VR-BS-GRA-000-006 - Papa Blanca Chilena (was "Papa")
}

foreach ($row in $categoryFixes.Keys) {
    $ws.Cells.Item($row, 3).Value = $categoryFixes[$row]
}

# --- column I (path_image): drop the leading "/src/" on every data row ----
for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 9)
    $path = $cell.Value()
    if ($path -like "/src/*") {
        $cell.Value = $path.Substring(5)
    }
}
